$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for the new "Save" column, reusing the same style as the
# other header cells (e.g. G1) so the bold/border/centered formatting
# matches exactly instead of minting a new style entry.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the Save values for each data row (0 = not saved, 1 = saved)
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
